# Pass client timezone to excel template - adjust timezone formatting
# in the events export template.
#
# - Replace the "Period:" value formula with one that formats the
#   from/to range using Joda-Time's toString pattern.
# - Replace the per-event "Time" column formula (event.serverTime) with
#   one that converts the server time into the client's timezone using
#   org.joda.time.DateTime before formatting it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Period:" range cell (row 6, column B) first so the new
# shared string for the from/to range is allocated before the new
# shared string used by the per-event time column (matches the order
# the strings appear in the saved workbook).
$ws.Range("B6").Value = '${from.toString("YYYY.MM.dd HH:mm:ss")+" - "+to.toString("YYYY.MM.dd HH:mm:ss")}'

# Header row (row 8) stays the same text, rewritten here only because
# the shared-string table is being compacted (the old "Period" format
# string and "event.serverTime" string are no longer used).
$ws.Range("A8").Value = "Time"
$ws.Range("B8").Value = "Type"
$ws.Range("C8").Value = "Geofence Name"
$ws.Range("D8").Value = "Attributes"

# Template data row (row 9): the server time column now converts the
# event's server time into the client timezone before formatting it.
$ws.Range("A9").Value = '${new("org.joda.time.DateTime", event.serverTime, timezone).toString("YYYY.MM.dd HH:mm:ss")}'
$ws.Range("B9").Value = '${event.type}'
$ws.Range("C9").Value = '${geofenceNames[event.geofenceId]}'
$ws.Range("D9").Value = '${event.attributes.toString().replaceAll(",", " ").replaceAll(bracketsRegex, "")}'
